$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.996.60'
$ws.Range("E2").Value = '  +0.39%  '
$ws.Range("D3").Value = '1.560.81'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.49'
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("E6").Value = '  +0.81%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.13'
$ws.Range("E8").Value = '  +2.08%  '
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("E10").Value = '  +1.96%  '
$ws.Range("E11").Value = '  +0.15%  '
$ws.Range("E12").Value = '  +0.58%  '
$ws.Range("D13").Value = '1.544.78'
$ws.Range("E13").Value = '  -0.44%  '
$ws.Range("E14").Value = '  +0.90%  '
$ws.Range("E15").Value = '  +1.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.08'
$ws.Range("E16").Value = '  +0.72%  '
$ws.Range("D17").Value = '26.995.44'
$ws.Range("E17").Value = '  +0.39%  '
$ws.Range("E18").Value = '  +2.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '217.43'
$ws.Range("E19").Value = '  +0.33%  '
$ws.Range("E20").Value = '  +2.71%  '
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("E22").Value = '  +1.58%  '
$ws.Range("E23").Value = '  +0.46%  '
$ws.Range("E24").Value = '  -1.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.63'
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("E26").Value = '  +0.34%  '
$ws.Range("E27").Value = '  +1.27%  '
$ws.Range("E28").Value = '  +1.65%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("E30").Value = '  +0.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.11'
$ws.Range("E31").Value = '  +1.97%  '
$ws.Range("E32").Value = '  +0.87%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.14'
$ws.Range("E33").Value = '  +4.46%  '
$ws.Range("B34").Value = 'Maker'
$ws.Range("C34").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D34").Value = '1.424.91'
$ws.Range("E34").Value = '  -0.08%  '
$ws.Range("E35").Value = '  +3.14%  '
$ws.Range("E36").Value = '  +9.83%  '
$ws.Range("E37").Value = '  +1.29%  '
$ws.Range("E38").Value = '  +0.88%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.534'
$ws.Range("E39").Value = '  +2.48%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.810'
$ws.Range("E40").Value = '  +0.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.01'
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.70'
$ws.Range("E42").Value = '  +0.14%  '
$ws.Range("E43").Value = '  +2.70%  '
$ws.Range("E44").Value = '  +1.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.91'
$ws.Range("E45").Value = '  +2.27%  '
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("D47").Value = '1.696.93'
$ws.Range("E47").Value = '  +0.64%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.43'
$ws.Range("E48").Value = '  +1.49%  '
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("D50").Value = '0.0₆01000'
$ws.Range("E50").Value = '  -1.08%  '
$ws.Range("E51").Value = '  -0.23%  '
